$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "67.497.71"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.630.46"
$ws.Range("E3").Value = "  -1.61%  "
Set-TextValue $ws.Range("D5") "594.57"
$ws.Range("E5").Value = "  -0.45%  "
Set-TextValue $ws.Range("D6") "168.67"
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").Value = "2.630.14"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +1.92%  "
Set-TextValue $ws.Range("D13") "5.23"
$ws.Range("E13").Value = "  +0.05%  "
Set-TextValue $ws.Range("D14") "27.71"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "3.109.34"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "67.295.43"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "2.629.42"
$ws.Range("E18").Value = "  -1.56%  "
Set-TextValue $ws.Range("D19") "12.03"
$ws.Range("E19").Value = "  +2.75%  "
Set-TextValue $ws.Range("D20") "8.05"
$ws.Range("E20").Value = "  +4.63%  "
Set-TextValue $ws.Range("D21") "356.74"
$ws.Range("E21").Value = "  -1.84%  "
$ws.Range("E22").Value = "  -1.14%  "
Set-TextValue $ws.Range("D23") "4.68"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D24") "1.94"
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  +0.02%  "
Set-TextValue $ws.Range("D26") "10.35"
$ws.Range("E26").Value = "  +3.37%  "
Set-TextValue $ws.Range("D27") "69.58"
$ws.Range("D28").Value = "2.767.40"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  -1.38%  "
Set-TextValue $ws.Range("D31") "547.45"
$ws.Range("E31").Value = "  -1.58%  "
Set-TextValue $ws.Range("D32") "7.94"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  -2.75%  "
Set-TextValue $ws.Range("D34") "1.89"
$ws.Range("E34").Value = "  -1.84%  "
Set-TextValue $ws.Range("D35") "0.136"
$ws.Range("E35").Value = "  +4.60%  "
$ws.Range("E36").Value = "  +0.08%  "
Set-TextValue $ws.Range("D37") "1.51"
$ws.Range("E37").Value = "  -3.02%  "
Set-TextValue $ws.Range("D38") "156.42"
$ws.Range("E38").Value = "  +1.20%  "
Set-TextValue $ws.Range("D39") "19.05"
$ws.Range("E39").Value = "  -2.44%  "
Set-TextValue $ws.Range("D40") "0.367"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("E41").Value = "  -0.54%  "
Set-TextValue $ws.Range("D42") "5.23"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("D46").Value = "0.0₆0298"
$ws.Range("E46").Value = "  -0.14%  "
Set-TextValue $ws.Range("D47") "152.90"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  -1.74%  "
Set-TextValue $ws.Range("D49") "3.80"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("E50").Value = "  -0.92%  "
Set-TextValue $ws.Range("D51") "0.0771"
$ws.Range("E51").Value = "  -1.14%  "
